# Apply the two substantive changes from the commit:
#  1. Re-style the table on slide 16 with the built-in "Medium Style 2 -
#     Accent 1" table style (its brace-GUID StyleId).
#  2. Re-colour the presentation's live theme (the one the slide master /
#     every slide actually renders with) from the custom "Integral" palette
#     to the standard "Office" palette -- i.e. switch the deck's colour
#     design from Integral to Office Theme.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{C6AA83FF-B70A-47E6-ACF5-51310FA988AC}")

# --- 2. Theme colours -------------------------------------------------------
# Office theme colour values (hex RRGGBB), in the standard 12-slot theme
# colour order used by ThemeColorScheme.Colors(index).
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
